$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7170026666666667
$ws.Range("H2").Value = 2.151008
$ws.Range("I2").Value = 0.02953485643833859
$ws.Range("J2").Value = 0.02953485643833859
$ws.Range("M2").Value = 16.57637
$ws.Range("N2").Value = 49.72911
$ws.Range("O2").Value = 0.1853914334114506
$ws.Range("P2").Value = 0.1853914334114506
$ws.Range("Q2").Value = 11.88530149365333
$ws.Range("R2").Value = 106.96771344288
$ws.Range("S2").Value = 0.005475509370705003
$ws.Range("T2").Value = 0.005475509370705004
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7170026666666667
$ws.Range("H3").Value = 2.151008
$ws.Range("I3").Value = 0.02953485643833859
$ws.Range("J3").Value = 0.02953485643833859
$ws.Range("O3").Value = 0.5978024790674488
$ws.Range("P3").Value = 0.5978024790674489
$ws.Range("Q3").Value = 38.32465484854045
$ws.Range("R3").Value = 344.921893636864
$ws.Range("S3").Value = 0.01765601039774001
$ws.Range("T3").Value = 0.01765601039774001
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7170026666666667
$ws.Range("H4").Value = 2.151008
$ws.Range("I4").Value = 0.02953485643833859
$ws.Range("J4").Value = 0.02953485643833859
$ws.Range("O4").Value = 0.2168060875211005
$ws.Range("P4").Value = 0.2168060875211005
$ws.Range("Q4").Value = 13.89927068598045
$ws.Range("R4").Value = 125.093436173824
$ws.Range("S4").Value = 0.006403336669893575
$ws.Range("T4").Value = 0.006403336669893575
$ws.Range("I5").Value = 0.4970672037825566
$ws.Range("J5").Value = 0.4970672037825566
$ws.Range("M5").Value = 16.57637
$ws.Range("N5").Value = 49.72911
$ws.Range("O5").Value = 0.1853914334114506
$ws.Range("P5").Value = 0.1853914334114506
$ws.Range("Q5").Value = 200.0278414048467
$ws.Range("R5").Value = 1800.25057264362
$ws.Range("S5").Value = 0.09215200141106981
$ws.Range("T5").Value = 0.09215200141106981
$ws.Range("I6").Value = 0.4970672037825566
$ws.Range("J6").Value = 0.4970672037825566
$ws.Range("O6").Value = 0.5978024790674488
$ws.Range("P6").Value = 0.5978024790674489
$ws.Range("Q6").Value = 644.9981926022596
$ws.Range("R6").Value = 5804.983733420336
$ws.Range("S6").Value = 0.2971480066843371
$ws.Range("T6").Value = 0.2971480066843371
$ws.Range("I7").Value = 0.4970672037825566
$ws.Range("J7").Value = 0.4970672037825566
$ws.Range("O7").Value = 0.2168060875211005
$ws.Range("P7").Value = 0.2168060875211005
$ws.Range("S7").Value = 0.1077671956871497
$ws.Range("T7").Value = 0.1077671956871497
$ws.Range("I8").Value = 0.4733979397791048
$ws.Range("J8").Value = 0.4733979397791048
$ws.Range("M8").Value = 16.57637
$ws.Range("N8").Value = 49.72911
$ws.Range("O8").Value = 0.1853914334114506
$ws.Range("P8").Value = 0.1853914334114506
$ws.Range("Q8").Value = 190.5029486937133
$ws.Range("R8").Value = 1714.52653824342
$ws.Range("S8").Value = 0.08776392262967583
$ws.Range("T8").Value = 0.08776392262967583
$ws.Range("I9").Value = 0.4733979397791048
$ws.Range("J9").Value = 0.4733979397791048
$ws.Range("O9").Value = 0.5978024790674488
$ws.Range("P9").Value = 0.5978024790674489
$ws.Range("S9").Value = 0.2829984619853717
$ws.Range("T9").Value = 0.2829984619853718
$ws.Range("I10").Value = 0.4733979397791048
$ws.Range("J10").Value = 0.4733979397791048
$ws.Range("O10").Value = 0.2168060875211005
$ws.Range("P10").Value = 0.2168060875211005
$ws.Range("S10").Value = 0.1026355551640573
$ws.Range("T10").Value = 0.1026355551640573
